$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 22:49:18"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 22:49:23"
